$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency list: updated prices/volumes for top ranks,
# and shifted coin rows (OKB and Frax dropped out of top 50, Aave and Cronos added).
# NumberFormat is forced to Text before writing so numeric-looking strings
# (e.g. "1.003", "0.3840") are preserved verbatim as text, matching the source data;
# the original cell style is restored immediately after so no formatting changes.
$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.068.62'
$ws.Range('D2').Style = $origStyle
$origStyle = $ws.Range('E2').Style
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('E2').Style = $origStyle
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.875.27'
$ws.Range('D3').Style = $origStyle
$origStyle = $ws.Range('E3').Style
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.31%  '
$ws.Range('E3').Style = $origStyle
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = $origStyle
$origStyle = $ws.Range('E4').Style
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('E4').Style = $origStyle
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '313.65'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('E5').Style
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.37%  '
$ws.Range('E5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = $origStyle
$origStyle = $ws.Range('E6').Style
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E6').Style = $origStyle
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5076'
$ws.Range('D7').Style = $origStyle
$origStyle = $ws.Range('E7').Style
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E7').Style = $origStyle
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3840'
$ws.Range('D8').Style = $origStyle
$origStyle = $ws.Range('E8').Style
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('E8').Style = $origStyle
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08977'
$ws.Range('D9').Style = $origStyle
$origStyle = $ws.Range('E9').Style
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -2.62%  '
$ws.Range('E9').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.123'
$ws.Range('D10').Style = $origStyle
$origStyle = $ws.Range('E10').Style
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('E10').Style = $origStyle
$origStyle = $ws.Range('B11').Style
$ws.Range('B11').NumberFormat = "@"
$ws.Range('B11').Value = 'Polkadot'
$ws.Range('B11').Style = $origStyle
$origStyle = $ws.Range('C11').Style
$ws.Range('C11').NumberFormat = "@"
$ws.Range('C11').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C11').Style = $origStyle
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.337'
$ws.Range('D11').Style = $origStyle
$origStyle = $ws.Range('E11').Style
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('E11').Style = $origStyle
$origStyle = $ws.Range('B12').Style
$ws.Range('B12').NumberFormat = "@"
$ws.Range('B12').Value = 'Solana'
$ws.Range('B12').Style = $origStyle
$origStyle = $ws.Range('C12').Style
$ws.Range('C12').NumberFormat = "@"
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('C12').Style = $origStyle
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '20.72'
$ws.Range('D12').Style = $origStyle
$origStyle = $ws.Range('E12').Style
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('E12').Style = $origStyle
$origStyle = $ws.Range('B13').Style
$ws.Range('B13').NumberFormat = "@"
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('B13').Style = $origStyle
$origStyle = $ws.Range('C13').Style
$ws.Range('C13').NumberFormat = "@"
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C13').Style = $origStyle
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.877.23'
$ws.Range('D13').Style = $origStyle
$origStyle = $ws.Range('E13').Style
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('E13').Style = $origStyle
$origStyle = $ws.Range('B14').Style
$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('B14').Style = $origStyle
$origStyle = $ws.Range('C14').Style
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C14').Style = $origStyle
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.206'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('E14').Style
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('E14').Style = $origStyle
$origStyle = $ws.Range('B15').Style
$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'BinanceUSD'
$ws.Range('B15').Style = $origStyle
$origStyle = $ws.Range('C15').Style
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('C15').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.003'
$ws.Range('D15').Style = $origStyle
$origStyle = $ws.Range('E15').Style
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('E15').Style = $origStyle
$origStyle = $ws.Range('B16').Style
$ws.Range('B16').NumberFormat = "@"
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('B16').Style = $origStyle
$origStyle = $ws.Range('C16').Style
$ws.Range('C16').NumberFormat = "@"
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('C16').Style = $origStyle
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001106'
$ws.Range('D16').Style = $origStyle
$origStyle = $ws.Range('E16').Style
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('E16').Style = $origStyle
$origStyle = $ws.Range('B17').Style
$ws.Range('B17').NumberFormat = "@"
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('B17').Style = $origStyle
$origStyle = $ws.Range('C17').Style
$ws.Range('C17').NumberFormat = "@"
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('C17').Style = $origStyle
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '91.14'
$ws.Range('D17').Style = $origStyle
$origStyle = $ws.Range('E17').Style
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('E17').Style = $origStyle
$origStyle = $ws.Range('B18').Style
$ws.Range('B18').NumberFormat = "@"
$ws.Range('B18').Value = 'TRON'
$ws.Range('B18').Style = $origStyle
$origStyle = $ws.Range('C18').Style
$ws.Range('C18').NumberFormat = "@"
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C18').Style = $origStyle
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06599'
$ws.Range('D18').Style = $origStyle
$origStyle = $ws.Range('E18').Style
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('E18').Style = $origStyle
$origStyle = $ws.Range('B19').Style
$ws.Range('B19').NumberFormat = "@"
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('B19').Style = $origStyle
$origStyle = $ws.Range('C19').Style
$ws.Range('C19').NumberFormat = "@"
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C19').Style = $origStyle
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '18.14'
$ws.Range('D19').Style = $origStyle
$origStyle = $ws.Range('E19').Style
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +2.14%  '
$ws.Range('E19').Style = $origStyle
$origStyle = $ws.Range('B20').Style
$ws.Range('B20').NumberFormat = "@"
$ws.Range('B20').Value = 'Dai'
$ws.Range('B20').Style = $origStyle
$origStyle = $ws.Range('C20').Style
$ws.Range('C20').NumberFormat = "@"
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('C20').Style = $origStyle
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = $origStyle
$origStyle = $ws.Range('E20').Style
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E20').Style = $origStyle
$origStyle = $ws.Range('B21').Style
$ws.Range('B21').NumberFormat = "@"
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('B21').Style = $origStyle
$origStyle = $ws.Range('C21').Style
$ws.Range('C21').NumberFormat = "@"
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C21').Style = $origStyle
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.115'
$ws.Range('D21').Style = $origStyle
$origStyle = $ws.Range('E21').Style
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.57%  '
$ws.Range('E21').Style = $origStyle
$origStyle = $ws.Range('B22').Style
$ws.Range('B22').NumberFormat = "@"
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('B22').Style = $origStyle
$origStyle = $ws.Range('C22').Style
$ws.Range('C22').NumberFormat = "@"
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C22').Style = $origStyle
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '28.097.51'
$ws.Range('D22').Style = $origStyle
$origStyle = $ws.Range('E22').Style
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E22').Style = $origStyle
$origStyle = $ws.Range('B23').Style
$ws.Range('B23').NumberFormat = "@"
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('B23').Style = $origStyle
$origStyle = $ws.Range('C23').Style
$ws.Range('C23').NumberFormat = "@"
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C23').Style = $origStyle
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '11.39'
$ws.Range('D23').Style = $origStyle
$origStyle = $ws.Range('E23').Style
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.50%  '
$ws.Range('E23').Style = $origStyle
$origStyle = $ws.Range('B24').Style
$ws.Range('B24').NumberFormat = "@"
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('B24').Style = $origStyle
$origStyle = $ws.Range('C24').Style
$ws.Range('C24').NumberFormat = "@"
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C24').Style = $origStyle
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.271'
$ws.Range('D24').Style = $origStyle
$origStyle = $ws.Range('E24').Style
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('E24').Style = $origStyle
$origStyle = $ws.Range('B25').Style
$ws.Range('B25').NumberFormat = "@"
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('B25').Style = $origStyle
$origStyle = $ws.Range('C25').Style
$ws.Range('C25').NumberFormat = "@"
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C25').Style = $origStyle
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.094.00'
$ws.Range('D25').Style = $origStyle
$origStyle = $ws.Range('E25').Style
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -1.03%  '
$ws.Range('E25').Style = $origStyle
$origStyle = $ws.Range('B26').Style
$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('B26').Style = $origStyle
$origStyle = $ws.Range('C26').Style
$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C26').Style = $origStyle
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.536'
$ws.Range('D26').Style = $origStyle
$origStyle = $ws.Range('E26').Style
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -2.53%  '
$ws.Range('E26').Style = $origStyle
$origStyle = $ws.Range('B27').Style
$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('B27').Style = $origStyle
$origStyle = $ws.Range('C27').Style
$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C27').Style = $origStyle
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.75'
$ws.Range('D27').Style = $origStyle
$origStyle = $ws.Range('E27').Style
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('E27').Style = $origStyle
$origStyle = $ws.Range('B28').Style
$ws.Range('B28').NumberFormat = "@"
$ws.Range('B28').Value = 'Monero'
$ws.Range('B28').Style = $origStyle
$origStyle = $ws.Range('C28').Style
$ws.Range('C28').NumberFormat = "@"
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C28').Style = $origStyle
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '156.96'
$ws.Range('D28').Style = $origStyle
$origStyle = $ws.Range('E28').Style
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('E28').Style = $origStyle
$origStyle = $ws.Range('B29').Style
$ws.Range('B29').NumberFormat = "@"
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('B29').Style = $origStyle
$origStyle = $ws.Range('C29').Style
$ws.Range('C29').NumberFormat = "@"
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C29').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '126.74'
$ws.Range('D29').Style = $origStyle
$origStyle = $ws.Range('E29').Style
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('E29').Style = $origStyle
$origStyle = $ws.Range('B30').Style
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'Stellar'
$ws.Range('B30').Style = $origStyle
$origStyle = $ws.Range('C30').Style
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C30').Style = $origStyle
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.1052'
$ws.Range('D30').Style = $origStyle
$origStyle = $ws.Range('E30').Style
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('E30').Style = $origStyle
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.059'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('E31').Style
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.38%  '
$ws.Range('E31').Style = $origStyle
$origStyle = $ws.Range('B32').Style
$ws.Range('B32').NumberFormat = "@"
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('B32').Style = $origStyle
$origStyle = $ws.Range('C32').Style
$ws.Range('C32').NumberFormat = "@"
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C32').Style = $origStyle
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.602'
$ws.Range('D32').Style = $origStyle
$origStyle = $ws.Range('E32').Style
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('E32').Style = $origStyle
$origStyle = $ws.Range('B33').Style
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('B33').Style = $origStyle
$origStyle = $ws.Range('C33').Style
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C33').Style = $origStyle
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.601'
$ws.Range('D33').Style = $origStyle
$origStyle = $ws.Range('E33').Style
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('E33').Style = $origStyle
$origStyle = $ws.Range('B34').Style
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'FraxShare'
$ws.Range('B34').Style = $origStyle
$origStyle = $ws.Range('C34').Style
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C34').Style = $origStyle
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '9.613'
$ws.Range('D34').Style = $origStyle
$origStyle = $ws.Range('E34').Style
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('E34').Style = $origStyle
$origStyle = $ws.Range('B35').Style
$ws.Range('B35').NumberFormat = "@"
$ws.Range('B35').Value = 'Hedera'
$ws.Range('B35').Style = $origStyle
$origStyle = $ws.Range('C35').Style
$ws.Range('C35').NumberFormat = "@"
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C35').Style = $origStyle
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.06574'
$ws.Range('D35').Style = $origStyle
$origStyle = $ws.Range('E35').Style
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('E35').Style = $origStyle
$origStyle = $ws.Range('B36').Style
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'VeChain'
$ws.Range('B36').Style = $origStyle
$origStyle = $ws.Range('C36').Style
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C36').Style = $origStyle
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02417'
$ws.Range('D36').Style = $origStyle
$origStyle = $ws.Range('E36').Style
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +0.54%  '
$ws.Range('E36').Style = $origStyle
$origStyle = $ws.Range('B37').Style
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Algorand'
$ws.Range('B37').Style = $origStyle
$origStyle = $ws.Range('C37').Style
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C37').Style = $origStyle
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2178'
$ws.Range('D37').Style = $origStyle
$origStyle = $ws.Range('E37').Style
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.46%  '
$ws.Range('E37').Style = $origStyle
$origStyle = $ws.Range('B38').Style
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('B38').Style = $origStyle
$origStyle = $ws.Range('C38').Style
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C38').Style = $origStyle
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.275'
$ws.Range('D38').Style = $origStyle
$origStyle = $ws.Range('E38').Style
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +1.36%  '
$ws.Range('E38').Style = $origStyle
$origStyle = $ws.Range('B39').Style
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('B39').Style = $origStyle
$origStyle = $ws.Range('C39').Style
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C39').Style = $origStyle
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.206'
$ws.Range('D39').Style = $origStyle
$origStyle = $ws.Range('E39').Style
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('E39').Style = $origStyle
$origStyle = $ws.Range('B40').Style
$ws.Range('B40').NumberFormat = "@"
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('B40').Style = $origStyle
$origStyle = $ws.Range('C40').Style
$ws.Range('C40').NumberFormat = "@"
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C40').Style = $origStyle
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.6396'
$ws.Range('D40').Style = $origStyle
$origStyle = $ws.Range('E40').Style
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.62%  '
$ws.Range('E40').Style = $origStyle
$origStyle = $ws.Range('B41').Style
$ws.Range('B41').NumberFormat = "@"
$ws.Range('B41').Value = 'Aptos'
$ws.Range('B41').Style = $origStyle
$origStyle = $ws.Range('C41').Style
$ws.Range('C41').NumberFormat = "@"
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C41').Style = $origStyle
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '11.47'
$ws.Range('D41').Style = $origStyle
$origStyle = $ws.Range('E41').Style
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.72%  '
$ws.Range('E41').Style = $origStyle
$origStyle = $ws.Range('B42').Style
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'InternetComputer(DFINITY)'
$ws.Range('B42').Style = $origStyle
$origStyle = $ws.Range('C42').Style
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C42').Style = $origStyle
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.917'
$ws.Range('D42').Style = $origStyle
$origStyle = $ws.Range('E42').Style
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -1.24%  '
$ws.Range('E42').Style = $origStyle
$origStyle = $ws.Range('B43').Style
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('B43').Style = $origStyle
$origStyle = $ws.Range('C43').Style
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C43').Style = $origStyle
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.24'
$ws.Range('D43').Style = $origStyle
$origStyle = $ws.Range('E43').Style
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('E43').Style = $origStyle
$origStyle = $ws.Range('B44').Style
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('B44').Style = $origStyle
$origStyle = $ws.Range('C44').Style
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C44').Style = $origStyle
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6026'
$ws.Range('D44').Style = $origStyle
$origStyle = $ws.Range('E44').Style
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.87%  '
$ws.Range('E44').Style = $origStyle
$origStyle = $ws.Range('B45').Style
$ws.Range('B45').NumberFormat = "@"
$ws.Range('B45').Value = 'PancakeSwap'
$ws.Range('B45').Style = $origStyle
$origStyle = $ws.Range('C45').Style
$ws.Range('C45').NumberFormat = "@"
$ws.Range('C45').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C45').Style = $origStyle
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.675'
$ws.Range('D45').Style = $origStyle
$origStyle = $ws.Range('E45').Style
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('E45').Style = $origStyle
$origStyle = $ws.Range('B46').Style
$ws.Range('B46').NumberFormat = "@"
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('B46').Style = $origStyle
$origStyle = $ws.Range('C46').Style
$ws.Range('C46').NumberFormat = "@"
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C46').Style = $origStyle
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.276'
$ws.Range('D46').Style = $origStyle
$origStyle = $ws.Range('E46').Style
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E46').Style = $origStyle
$origStyle = $ws.Range('B47').Style
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'EOS'
$ws.Range('B47').Style = $origStyle
$origStyle = $ws.Range('C47').Style
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('C47').Style = $origStyle
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.239'
$ws.Range('D47').Style = $origStyle
$origStyle = $ws.Range('E47').Style
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +5.23%  '
$ws.Range('E47').Style = $origStyle
$origStyle = $ws.Range('B48').Style
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('B48').Style = $origStyle
$origStyle = $ws.Range('C48').Style
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C48').Style = $origStyle
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.995'
$ws.Range('D48').Style = $origStyle
$origStyle = $ws.Range('E48').Style
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.56%  '
$ws.Range('E48').Style = $origStyle
$origStyle = $ws.Range('B49').Style
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'Quant'
$ws.Range('B49').Style = $origStyle
$origStyle = $ws.Range('C49').Style
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C49').Style = $origStyle
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '121.28'
$ws.Range('D49').Style = $origStyle
$origStyle = $ws.Range('E49').Style
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.93%  '
$ws.Range('E49').Style = $origStyle
$origStyle = $ws.Range('B50').Style
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'Aave'
$ws.Range('B50').Style = $origStyle
$origStyle = $ws.Range('C50').Style
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('C50').Style = $origStyle
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '79.65'
$ws.Range('D50').Style = $origStyle
$origStyle = $ws.Range('E50').Style
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +2.16%  '
$ws.Range('E50').Style = $origStyle
$origStyle = $ws.Range('B51').Style
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'Cronos'
$ws.Range('B51').Style = $origStyle
$origStyle = $ws.Range('C51').Style
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C51').Style = $origStyle
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06936'
$ws.Range('D51').Style = $origStyle
$origStyle = $ws.Range('E51').Style
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.50%  '
$ws.Range('E51').Style = $origStyle
